$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028146505764195
$ws.Cells.Item(2, 4).Value = 1.036210191072655
$ws.Cells.Item(2, 5).Value = 1.028150765179679
$ws.Cells.Item(2, 6).Value = 1.044933865562794
$ws.Cells.Item(2, 9).Value = 1.032497490467126
$ws.Cells.Item(2, 10).Value = 1.03330044293342
$ws.Cells.Item(2, 11).Value = 1.039004584774202
$ws.Cells.Item(2, 12).Value = 1.030968418151578
$ws.Cells.Item(2, 13).Value = 1.047703521851757
$ws.Cells.Item(2, 14).Value = 1.01510907989309
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028988062831014
$ws.Cells.Item(3, 4).Value = 1.03696019127108
$ws.Cells.Item(3, 5).Value = 1.028862626433227
$ws.Cells.Item(3, 6).Value = 1.045781930100127
$ws.Cells.Item(3, 9).Value = 1.032587526028147
$ws.Cells.Item(3, 10).Value = 1.033783039078247
$ws.Cells.Item(3, 11).Value = 1.039564297504068
$ws.Cells.Item(3, 12).Value = 1.031488413737043
$ws.Cells.Item(3, 13).Value = 1.048362829051356
$ws.Cells.Item(3, 14).Value = 1.015270328852014
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.0295332386819
$ws.Cells.Item(4, 4).Value = 1.037446386082882
$ws.Cells.Item(4, 5).Value = 1.029324180186136
$ws.Cells.Item(4, 6).Value = 1.04633168560445
$ws.Cells.Item(4, 9).Value = 1.032644498079739
$ws.Cells.Item(4, 10).Value = 1.0340952938377
$ws.Cells.Item(4, 11).Value = 1.039926698513201
$ws.Cells.Item(4, 12).Value = 1.031825145411283
$ws.Cells.Item(4, 13).Value = 1.048789800276794
$ws.Cells.Item(4, 14).Value = 1.015374620407654
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029762580286592
$ws.Cells.Item(5, 4).Value = 1.037650994660433
$ws.Cells.Item(5, 5).Value = 1.029518438638771
$ws.Cells.Item(5, 6).Value = 1.046563040477852
$ws.Cells.Item(5, 9).Value = 1.032668140369346
$ws.Cells.Item(5, 10).Value = 1.034226560296859
$ws.Cells.Item(5, 11).Value = 1.040079105344855
$ws.Cells.Item(5, 12).Value = 1.031966768181086
$ws.Cells.Item(5, 13).Value = 1.048969382362682
$ws.Cells.Item(5, 14).Value = 1.015418452770038
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.029801096480031
$ws.Cells.Item(6, 4).Value = 1.037685361734529
$ws.Cells.Item(6, 5).Value = 1.029551068408442
$ws.Cells.Item(6, 6).Value = 1.046601899857655
$ws.Cells.Item(6, 9).Value = 1.032672091888855
$ws.Cells.Item(6, 10).Value = 1.03424860015546
$ws.Cells.Item(6, 11).Value = 1.040104698213581
$ws.Cells.Item(6, 12).Value = 1.031990550790017
$ws.Cells.Item(6, 13).Value = 1.048999539825381
$ws.Cells.Item(6, 14).Value = 1.015425811710094
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.029536302569598
$ws.Cells.Item(7, 4).Value = 1.037449119240957
$ws.Cells.Item(7, 5).Value = 1.029326775010035
$ws.Cells.Item(7, 6).Value = 1.046334776048256
$ws.Cells.Item(7, 9).Value = 1.032644815203208
$ws.Cells.Item(7, 10).Value = 1.034097047849593
$ws.Cells.Item(7, 11).Value = 1.039928734772597
$ws.Cells.Item(7, 12).Value = 1.031827037544642
$ws.Cells.Item(7, 13).Value = 1.04879219953591
$ws.Cells.Item(7, 14).Value = 1.015375206144829
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028430782022922
$ws.Cells.Item(8, 4).Value = 1.036463470805937
$ws.Cells.Item(8, 5).Value = 1.028391147951806
$ws.Cells.Item(8, 6).Value = 1.04522026447443
$ws.Cells.Item(8, 9).Value = 1.032528184489035
$ws.Cells.Item(8, 10).Value = 1.033463541331342
$ws.Cells.Item(8, 11).Value = 1.039193693912428
$ws.Cells.Item(8, 12).Value = 1.031144098460199
$ws.Cells.Item(8, 13).Value = 1.047926263203138
$ws.Cells.Item(8, 14).Value = 1.015163584237835
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026487628706017
$ws.Cells.Item(9, 4).Value = 1.034733566229871
$ws.Cells.Item(9, 5).Value = 1.026749672649462
$ws.Cells.Item(9, 6).Value = 1.043264107537928
$ws.Cells.Item(9, 9).Value = 1.032312838580192
$ws.Cells.Item(9, 10).Value = 1.032347148401942
$ws.Cells.Item(9, 11).Value = 1.03790028081943
$ws.Cells.Item(9, 12).Value = 1.029942724434432
$ws.Cells.Item(9, 13).Value = 1.046403169537221
$ws.Cells.Item(9, 14).Value = 1.014790338458625
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025195596822348
$ws.Cells.Item(10, 4).Value = 1.033585070293824
$ws.Cells.Item(10, 5).Value = 1.025660319844918
$ws.Cells.Item(10, 6).Value = 1.041965335672788
$ws.Cells.Item(10, 9).Value = 1.032162706183009
$ws.Cells.Item(10, 10).Value = 1.031602915709465
$ws.Cells.Item(10, 11).Value = 1.037039319374615
$ws.Cells.Item(10, 12).Value = 1.029143272256776
$ws.Cells.Item(10, 13).Value = 1.045389755328862
$ws.Cells.Item(10, 14).Value = 1.014541307494185
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024636960868446
$ws.Cells.Item(11, 4).Value = 1.033088914783689
$ws.Cells.Item(11, 5).Value = 1.025189818802213
$ws.Cells.Item(11, 6).Value = 1.041404243073465
$ws.Cells.Item(11, 9).Value = 1.03209614748237
$ws.Cells.Item(11, 10).Value = 1.031280678271749
$ws.Cells.Item(11, 11).Value = 1.036666843871825
$ws.Cells.Item(11, 12).Value = 1.028797465634837
$ws.Cells.Item(11, 13).Value = 1.044951427145161
$ws.Cells.Item(11, 14).Value = 1.014433432867183
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024429583823757
$ws.Cells.Item(12, 4).Value = 1.032904795349243
$ws.Cells.Item(12, 5).Value = 1.025015235334995
$ws.Cells.Item(12, 6).Value = 1.041196023267495
$ws.Cells.Item(12, 9).Value = 1.032071192336116
$ws.Cells.Item(12, 10).Value = 1.031160989313671
$ws.Cells.Item(12, 11).Value = 1.036528540441229
$ws.Cells.Item(12, 12).Value = 1.028669073399225
$ws.Cells.Item(12, 13).Value = 1.044788687409395
$ws.Cells.Item(12, 14).Value = 1.014393357519652
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024474061210421
$ws.Cells.Item(13, 4).Value = 1.032944281651502
$ws.Cells.Item(13, 5).Value = 1.025052675834534
$ws.Cells.Item(13, 6).Value = 1.041240678267339
$ws.Cells.Item(13, 9).Value = 1.032076555800539
$ws.Cells.Item(13, 10).Value = 1.031186662782119
$ws.Cells.Item(13, 11).Value = 1.03655820468467
$ws.Cells.Item(13, 12).Value = 1.028696611424144
$ws.Cells.Item(13, 13).Value = 1.044823592219536
$ws.Cells.Item(13, 14).Value = 1.014401954077243
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024619816453089
$ws.Cells.Item(14, 4).Value = 1.033073691836581
$ws.Cells.Item(14, 5).Value = 1.025175383963435
$ws.Cells.Item(14, 6).Value = 1.041387027570303
$ws.Cells.Item(14, 9).Value = 1.032094089416925
$ws.Cells.Item(14, 10).Value = 1.031270784648193
$ws.Cells.Item(14, 11).Value = 1.036655410627088
$ws.Cells.Item(14, 12).Value = 1.028786851547184
$ws.Cells.Item(14, 13).Value = 1.04493797349058
$ws.Cells.Item(14, 14).Value = 1.01443012034415
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024709637711711
$ws.Cells.Item(15, 4).Value = 1.033153448944463
$ws.Cells.Item(15, 5).Value = 1.025251012607897
$ws.Cells.Item(15, 6).Value = 1.041477224101962
$ws.Cells.Item(15, 9).Value = 1.032104861696472
$ws.Cells.Item(15, 10).Value = 1.031322615525752
$ws.Cells.Item(15, 11).Value = 1.036715309160525
$ws.Cells.Item(15, 12).Value = 1.028842458894465
$ws.Cells.Item(15, 13).Value = 1.04500845753904
$ws.Cells.Item(15, 14).Value = 1.014447473742542
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025232689119727
$ws.Cells.Item(16, 4).Value = 1.033618022910909
$ws.Cells.Item(16, 5).Value = 1.025691570808016
$ws.Cells.Item(16, 6).Value = 1.042002600770643
$ws.Cells.Item(16, 9).Value = 1.03216709085513
$ws.Cells.Item(16, 10).Value = 1.031624302087092
$ws.Cells.Item(16, 11).Value = 1.037064046354155
$ws.Cells.Item(16, 12).Value = 1.029166230041391
$ws.Cells.Item(16, 13).Value = 1.045418856131452
$ws.Cells.Item(16, 14).Value = 1.014548465922596
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02556100680022
$ws.Cells.Item(17, 4).Value = 1.033909747344482
$ws.Cells.Item(17, 5).Value = 1.025968242846734
$ws.Cells.Item(17, 6).Value = 1.042332500923953
$ws.Cells.Item(17, 9).Value = 1.032205710936871
$ws.Cells.Item(17, 10).Value = 1.031813548454886
$ws.Cells.Item(17, 11).Value = 1.037282888454033
$ws.Cells.Item(17, 12).Value = 1.029369420912666
$ws.Cells.Item(17, 13).Value = 1.045676419887996
$ws.Cells.Item(17, 14).Value = 1.014611804595261
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02575258818401
$ws.Cells.Item(18, 4).Value = 1.034080016178905
$ws.Cells.Item(18, 5).Value = 1.026129736226913
$ws.Cells.Item(18, 6).Value = 1.042525049871188
$ws.Cells.Item(18, 9).Value = 1.032228087744415
$ws.Cells.Item(18, 10).Value = 1.031923934521779
$ws.Cells.Item(18, 11).Value = 1.037410566748487
$ws.Cells.Item(18, 12).Value = 1.029487973488177
$ws.Cells.Item(18, 13).Value = 1.045826699328537
$ws.Cells.Item(18, 14).Value = 1.014648744830712
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025817925873955
$ws.Cells.Item(19, 4).Value = 1.03413809220485
$ws.Cells.Item(19, 5).Value = 1.026184820788689
$ws.Cells.Item(19, 6).Value = 1.042590725017837
$ws.Cells.Item(19, 9).Value = 1.032235692257566
$ws.Cells.Item(19, 10).Value = 1.031961573593232
$ws.Cells.Item(19, 11).Value = 1.037454107027597
$ws.Cells.Item(19, 12).Value = 1.029528402707953
$ws.Cells.Item(19, 13).Value = 1.045877948629844
$ws.Cells.Item(19, 14).Value = 1.014661339779169
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02552577321253
$ws.Cells.Item(20, 4).Value = 1.033878436588052
$ws.Cells.Item(20, 5).Value = 1.025938546619712
$ws.Cells.Item(20, 6).Value = 1.042297092942938
$ws.Cells.Item(20, 9).Value = 1.032201582840153
$ws.Cells.Item(20, 10).Value = 1.031793243908367
$ws.Cells.Item(20, 11).Value = 1.037259405509095
$ws.Cells.Item(20, 12).Value = 1.02934761684753
$ws.Cells.Item(20, 13).Value = 1.045648780875195
$ws.Cells.Item(20, 14).Value = 1.014605009379881
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024576891696824
$ws.Cells.Item(21, 4).Value = 1.033035578921446
$ws.Cells.Item(21, 5).Value = 1.025139244450475
$ws.Cells.Item(21, 6).Value = 1.041343925949842
$ws.Cells.Item(21, 9).Value = 1.032088932610823
$ws.Cells.Item(21, 10).Value = 1.031246012723859
$ws.Cells.Item(21, 11).Value = 1.036626784492174
$ws.Cells.Item(21, 12).Value = 1.028760276529867
$ws.Cells.Item(21, 13).Value = 1.044904288976567
$ws.Cells.Item(21, 14).Value = 1.014421826239739
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023981017550297
$ws.Cells.Item(22, 4).Value = 1.032506653217384
$ws.Cells.Item(22, 5).Value = 1.024637742892473
$ws.Cells.Item(22, 6).Value = 1.040745761061337
$ws.Cells.Item(22, 9).Value = 1.032016761352814
$ws.Cells.Item(22, 10).Value = 1.030901972725047
$ws.Cells.Item(22, 11).Value = 1.03622932379378
$ws.Cells.Item(22, 12).Value = 1.028391315570527
$ws.Cells.Item(22, 13).Value = 1.044436631388355
$ws.Cells.Item(22, 14).Value = 1.014306617788937
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024296832399335
$ws.Cells.Item(23, 4).Value = 1.032786950140576
$ws.Cells.Item(23, 5).Value = 1.024903498059943
$ws.Cells.Item(23, 6).Value = 1.041062751888515
$ws.Cells.Item(23, 9).Value = 1.032055147842906
$ws.Cells.Item(23, 10).Value = 1.031084351931012
$ws.Cells.Item(23, 11).Value = 1.03643999693176
$ws.Cells.Item(23, 12).Value = 1.028586877639803
$ws.Cells.Item(23, 13).Value = 1.044684503809028
$ws.Cells.Item(23, 14).Value = 1.014367695022637
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025541693509567
$ws.Cells.Item(24, 4).Value = 1.033892584228685
$ws.Cells.Item(24, 5).Value = 1.025951964710553
$ws.Cells.Item(24, 6).Value = 1.042313091903269
$ws.Cells.Item(24, 9).Value = 1.032203448611787
$ws.Cells.Item(24, 10).Value = 1.03180241865335
$ws.Cells.Item(24, 11).Value = 1.037270016344267
$ws.Cells.Item(24, 12).Value = 1.029357469058743
$ws.Cells.Item(24, 13).Value = 1.0456612696108
$ws.Cells.Item(24, 14).Value = 1.01460807985807
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026989387756291
$ws.Cells.Item(25, 4).Value = 1.035179954797975
$ws.Cells.Item(25, 5).Value = 1.027173166480893
$ws.Cells.Item(25, 6).Value = 1.043768889758629
$ws.Cells.Item(25, 9).Value = 1.032369671262171
$ws.Cells.Item(25, 10).Value = 1.032635763335673
$ws.Cells.Item(25, 11).Value = 1.038234433983243
$ws.Cells.Item(25, 12).Value = 1.030253056570069
$ws.Cells.Item(25, 13).Value = 1.046796583898486
$ws.Cells.Item(25, 14).Value = 1.01488686887458
